$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("382:385").Insert()

$data = @(
    @(3,"Femacal de La Calera","Coquimbo",44588,5,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Lapins","Primera",210,5500,6000,5738,"`$/bandeja 10 kilos","Región de O'Higgins",574,10),
    @(3,"Femacal de La Calera","Coquimbo",44588,5,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Lapins","Segunda",80,4500,4500,4500,"`$/bandeja 10 kilos","Región de O'Higgins",450,10),
    @(3,"Femacal de La Calera","Coquimbo",44588,5,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Santina","Primera",210,5500,6000,5762,"`$/bandeja 10 kilos","Región de O'Higgins",576,10),
    @(3,"Femacal de La Calera","Coquimbo",44588,5,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Santina","Segunda",120,4500,4500,4500,"`$/bandeja 10 kilos","Región de O'Higgins",450,10)
)

$startRow = 382
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Host "D382:" $ws.Cells.Item(382,4).Value() "K382:" $ws.Cells.Item(382,11).Value() "M382:" $ws.Cells.Item(382,13).Value()
Write-Host "D383:" $ws.Cells.Item(383,4).Value() "K383:" $ws.Cells.Item(383,11).Value() "L383:" $ws.Cells.Item(383,12).Value()
Write-Host "D384:" $ws.Cells.Item(384,4).Value() "K384:" $ws.Cells.Item(384,11).Value()
Write-Host "D385:" $ws.Cells.Item(385,4).Value() "K385:" $ws.Cells.Item(385,11).Value() "S385:" $ws.Cells.Item(385,19).Value()
Write-Host "D386:" $ws.Cells.Item(386,4).Value() "K386:" $ws.Cells.Item(386,11).Value()
Write-Host "D477:" $ws.Cells.Item(477,4).Value() "K477:" $ws.Cells.Item(477,11).Value()
